# Update the "Oguls Otistic Team" roster sheet:
#  - Add a new player row (Cason Wallace) at the top of the data table
#  - Re-sort/rewrite the remaining rows to match the refreshed roster order
#    (player / position / team values for existing players are unchanged,
#     only their row order shifts)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired data (in row order), starting at row 2
$data = @(
    @("Cason Wallace", "PG,SG", "Oklahoma City Thunder"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Scoot Henderson", "PG", "Portland Trail Blazers"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers")
)

$rowCount = $data.Count

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
